$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B, shifting existing B:E data to E:H
$ws.Range("B:D").EntireColumn.Insert()

# New header values for the three inserted columns
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the new columns for existing data rows (2-27) with "UN" (unchanged marker)
$ws.Range("B2:D27").Value = "UN"

# Add two new analyst rows at the bottom
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
